$wb = $excel.ActiveWorkbook

# --- Rename sheets ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Producto"
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Inventarios"

# --- Sheet1 "Producto": update existing rows 2-4 values ---
$ws1.Range("A2").Value = "NOMBRE"
$ws1.Range("B2").Value = "Cerveza Artesanal"

$ws1.Range("A3").Value = "DESCRIPCIÓN"
$ws1.Range("B3").Value = "Cerveza artesanal de alta calidad"

$ws1.Range("A4").Value = "MARCA"
$ws1.Range("B4").Value = "Artesanal"

# --- Sheet1 "Producto": add new rows 5-9, replicating the alternating style pattern ---
# Row5 mirrors style of row3 (fill/alignment s=3). B5 holds a numeric-looking
# string ("50") so force text type first (apostrophe) before the format paste.
$ws1.Range("B5").Value = "'50"
$ws1.Range("A3:B3").Copy()
$ws1.Range("A5:B5").PasteSpecial(-4122)
$ws1.Range("A5").Value = "CONTENIDO"

# Row6 mirrors style of row4 (s=2)
$ws1.Range("A4:B4").Copy()
$ws1.Range("A6:B6").PasteSpecial(-4122)
$ws1.Range("A6").Value = "UNIDAD DE MEDIDA"
$ws1.Range("B6").Value = "L"

# Row7 mirrors style of row5/row3 (s=3)
$ws1.Range("A5:B5").Copy()
$ws1.Range("A7:B7").PasteSpecial(-4122)
$ws1.Range("A7").Value = "TIPO"
$ws1.Range("B7").Value = "Alcohólico"

# Row8 mirrors style of row6/row4 (s=2). B8 holds a numeric-looking string
# ("1500") so force text type first (apostrophe) before the format paste.
$ws1.Range("B8").Value = "'1500"
$ws1.Range("A6:B6").Copy()
$ws1.Range("A8:B8").PasteSpecial(-4122)
$ws1.Range("A8").Value = "PRECIO"

# Row9 mirrors style of row7/row5 (s=3)
$ws1.Range("A7:B7").Copy()
$ws1.Range("A9:B9").PasteSpecial(-4122)
$ws1.Range("A9").Value = "CATEGORÍA"
$ws1.Range("B9").Value = "Cerveza"

$excel.CutCopyMode = 0

# --- Sheet2 "Inventarios": remove column E entirely ---
$ws2.Range("E1:E3").Delete()

# --- Sheet2 "Inventarios": update header row ---
$ws2.Range("A1").Value = "NOMBRE"
$ws2.Range("B1").Value = "CANTIDAD"
$ws2.Range("C1").Value = "MÁXIMO STOCK"
$ws2.Range("D1").Value = "FECHA DE ACTUALIZACIÓN"

# --- Sheet2 "Inventarios": update data row 2 ---
$ws2.Range("A2").Value = "asdasdasdsa"
$ws2.Range("B2").Value = 50
$ws2.Range("C2").Value = 569
$ws2.Range("D2").Value = "2024-10-07T20:29:04.937Z"

# --- Sheet2 "Inventarios": update data row 3 ---
$ws2.Range("A3").Value = "?test"
$ws2.Range("B3").Value = 5555
$ws2.Range("C3").Value = 20000
$ws2.Range("D3").Value = "2024-10-15T00:03:30.008Z"

$excel.CutCopyMode = 0
